$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A261").Value = "17-09-2021"
$ws.Range("B261").Value = 2167
$ws.Range("C261").Value = 2704
$ws.Range("D261").Value = 15639
$ws.Range("E261").Value = 2730
$ws.Range("F261").Value = 4727
$ws.Range("G261").Value = 9081

$ws.Range("A262").Value = "18-09-2021"
$ws.Range("B262").Value = 2167
$ws.Range("C262").Value = 2704
$ws.Range("D262").Value = 15639
$ws.Range("E262").Value = 2730
$ws.Range("F262").Value = 4727
$ws.Range("G262").Value = 9081

$ws.Range("A263").Value = "19-09-2021"
$ws.Range("B263").Value = 2167
$ws.Range("C263").Value = 2704
$ws.Range("D263").Value = 15639
$ws.Range("E263").Value = 2730
$ws.Range("F263").Value = 4727
$ws.Range("G263").Value = 9081

$ws.Range("A264").Value = "20-09-2021"
$ws.Range("B264").Value = 2171
$ws.Range("C264").Value = 2708
$ws.Range("D264").Value = 15664
$ws.Range("E264").Value = 2734
$ws.Range("F264").Value = 4735
$ws.Range("G264").Value = 9096

$ws.Range("A265").Value = "21-09-2021"
$ws.Range("B265").Value = 2150
$ws.Range("C265").Value = 2683
$ws.Range("D265").Value = 15516
$ws.Range("E265").Value = 2708
$ws.Range("F265").Value = 4690
$ws.Range("G265").Value = 9010

$ws.Range("A266").Value = "22-09-2021"
$ws.Range("B266").Value = 2158
$ws.Range("C266").Value = 2693
$ws.Range("D266").Value = 15574
$ws.Range("E266").Value = 2719
$ws.Range("F266").Value = 4707
$ws.Range("G266").Value = 9043

$ws.Range("A267").Value = "23-09-2021"
$ws.Range("B267").Value = 2158
$ws.Range("C267").Value = 2693
$ws.Range("D267").Value = 15576
$ws.Range("E267").Value = 2719
$ws.Range("F267").Value = 4708
$ws.Range("G267").Value = 9044

$ws.Range("A268").Value = "24-09-2021"
$ws.Range("B268").Value = 2152
$ws.Range("C268").Value = 2686
$ws.Range("D268").Value = 15532
$ws.Range("E268").Value = 2711
$ws.Range("F268").Value = 4695
$ws.Range("G268").Value = 9019

$ws.Range("A269").Value = "25-09-2021"
$ws.Range("B269").Value = 2152
$ws.Range("C269").Value = 2686
$ws.Range("D269").Value = 15532
$ws.Range("E269").Value = 2711
$ws.Range("F269").Value = 4695
$ws.Range("G269").Value = 9019

$ws.Range("A270").Value = "26-09-2021"
$ws.Range("B270").Value = 2152
$ws.Range("C270").Value = 2686
$ws.Range("D270").Value = 15532
$ws.Range("E270").Value = 2711
$ws.Range("F270").Value = 4695
$ws.Range("G270").Value = 9019

$ws.Range("A271").Value = "27-09-2021"
$ws.Range("B271").Value = 2148
$ws.Range("C271").Value = 2680
$ws.Range("D271").Value = 15498
$ws.Range("E271").Value = 2705
$ws.Range("F271").Value = 4684
$ws.Range("G271").Value = 8999

$ws.Range("A272").Value = "28-09-2021"
$ws.Range("B272").Value = 2130
$ws.Range("C272").Value = 2658
$ws.Range("D272").Value = 15371
$ws.Range("E272").Value = 2683
$ws.Range("F272").Value = 4646
$ws.Range("G272").Value = 8925

$ws.Range("A273").Value = "29-09-2021"
$ws.Range("B273").Value = 2122
$ws.Range("C273").Value = 2647
$ws.Range("D273").Value = 15310
$ws.Range("E273").Value = 2672
$ws.Range("F273").Value = 4628
$ws.Range("G273").Value = 8890

$ws.Range("A274").Value = "30-09-2021"
$ws.Range("B274").Value = 2109
$ws.Range("C274").Value = 2631
$ws.Range("D274").Value = 15216
$ws.Range("E274").Value = 2656
$ws.Range("F274").Value = 4599
$ws.Range("G274").Value = 8835
